$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the entire row 534 ("「ライブ」" / beIN SPORTS live-broadcast post),
# which shifts all subsequent rows up by one.
$ws.Rows(534).Delete()
